# Auto-generated edit script to update cryptos.xlsx price/volume data
# per commit "Updated cryptos list on Mon Oct 14 06:39:48 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.073.83'
$ws.Range("E2").Value = '  +1.85%  '

# Row 3
$ws.Range("D3").Value = '2.528.91'
$ws.Range("E3").Value = '  +2.67%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '580.68'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  +1.14%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '151.94'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  +3.81%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.537'
$c.Style = 'Normal'
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("E9").Value = '  +0.49%  '

# Row 10
$ws.Range("E10").Value = '  -0.94%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '5.28'
$c.Style = 'Normal'
$ws.Range("E11").Value = '  -0.21%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '0.353'
$c.Style = 'Normal'
$ws.Range("E12").Value = '  -1.33%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '29.54'
$c.Style = 'Normal'
$ws.Range("E13").Value = '  +1.95%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '0.0000179'
$c.Style = 'Normal'
$ws.Range("E14").Value = '  +0.47%  '

# Row 15
$ws.Range("D15").Value = '2.976.91'
$ws.Range("E15").Value = '  +2.30%  '

# Row 16
$ws.Range("D16").Value = '63.879.25'
$ws.Range("E16").Value = '  +1.54%  '

# Row 17
$ws.Range("D17").Value = '2.529.62'
$ws.Range("E17").Value = '  +2.81%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '7.82'
$c.Style = 'Normal'
$ws.Range("E18").Value = '  -2.15%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '10.93'
$c.Style = 'Normal'
$ws.Range("E19").Value = '  -0.87%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '4.25'
$c.Style = 'Normal'
$ws.Range("E20").Value = '  +2.74%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '327.43'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  +0.06%  '

# Row 22
$ws.Range("E22").Value = '  +0.83%  '

# Row 23
$ws.Range("E23").Value = '  +0.00%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '10.10'
$c.Style = 'Normal'
$ws.Range("E24").Value = '  -0.52%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '65.40'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  -0.43%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '660.40'
$c.Style = 'Normal'
$ws.Range("E26").Value = '  +0.65%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '0.0000103'
$c.Style = 'Normal'
$ws.Range("E27").Value = '  +4.19%  '

# Row 28
$ws.Range("D28").Value = '2.647.73'
$ws.Range("E28").Value = '  +2.45%  '

# Row 29
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '1.49'
$c.Style = 'Normal'
$ws.Range("E29").Value = '  +2.87%  '

# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '0.995'
$c.Style = 'Normal'
$ws.Range("E30").Value = '  -0.60%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '8.01'
$c.Style = 'Normal'
$ws.Range("E31").Value = '  +0.06%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '1.86'
$c.Style = 'Normal'
$ws.Range("E32").Value = '  +0.26%  '

# Row 33
$ws.Range("E33").Value = '  +0.97%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range("E34").Value = '  -0.10%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '1.53'
$c.Style = 'Normal'
$ws.Range("E35").Value = '  -1.07%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '4.80'
$c.Style = 'Normal'
$ws.Range("E36").Value = '  +1.02%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '5.53'
$c.Style = 'Normal'
$ws.Range("E37").Value = '  +1.82%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '0.371'
$c.Style = 'Normal'
$ws.Range("E38").Value = '  +0.60%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '18.87'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  +0.75%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '151.89'
$c.Style = 'Normal'
$ws.Range("E40").Value = '  +1.03%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '2.80'
$c.Style = 'Normal'
$ws.Range("E41").Value = '  +0.96%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '1.77'
$c.Style = 'Normal'
$ws.Range("E42").Value = '  +1.94%  '

# Row 43
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '160.72'
$c.Style = 'Normal'
$ws.Range("E43").Value = '  +4.75%  '

# Row 44
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range("E44").Value = '  +0.05%  '

# Row 45
$ws.Range("D45").Value = '0.0₆0302'
$ws.Range("E45").Value = '  -1.80%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '15.42'
$c.Style = 'Normal'
$ws.Range("E46").Value = '  +1.17%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '3.65'
$c.Style = 'Normal'
$ws.Range("E47").Value = '  +1.60%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '20.93'
$c.Style = 'Normal'
$ws.Range("E48").Value = '  +1.63%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '0.617'
$c.Style = 'Normal'
$ws.Range("E49").Value = '  +1.58%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '0.0517'
$c.Style = 'Normal'
$ws.Range("E50").Value = '  +1.10%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '0.0229'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  +1.21%  '

Write-Output "Updated cryptos sheet values"
